$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set all cell values first (text + numeric) ---
# Row 7
$ws.Cells.Item(7, 1).Value = 44282
$ws.Cells.Item(7, 2).Value = 44288
$ws.Cells.Item(7, 3).Value = 44288
$ws.Cells.Item(7, 4).Value = "PR-006"
$ws.Cells.Item(7, 5).Value = "Closed"
$ws.Cells.Item(7, 6).Value = "Software-app"
$ws.Cells.Item(7, 7).Value = "Missing MIDI"
$ws.Cells.Item(7, 8).Value = "No MIDI output observed in behavior of synthesizer."

# Row 8
$ws.Cells.Item(8, 1).Value = 44295
$ws.Cells.Item(8, 2).Value = 44303
$ws.Cells.Item(8, 3).Value = 44303
$ws.Cells.Item(8, 4).Value = "PR-007"
$ws.Cells.Item(8, 5).Value = "Closed"
$ws.Cells.Item(8, 6).Value = "Software-app"
$ws.Cells.Item(8, 7).Value = "Missing MIDI Data"
$ws.Cells.Item(8, 8).Value = "Missing MIDI note data at the synthesizer."

# Row 9
$ws.Cells.Item(9, 1).Value = 44317
$ws.Cells.Item(9, 4).Value = "PR-008"
$ws.Cells.Item(9, 5).Value = "Open"
$ws.Cells.Item(9, 6).Value = "Software-app"
$ws.Cells.Item(9, 7).Value = "All functions missing"
$ws.Cells.Item(9, 8).Value = "Systems appears to have crashed. Recovered after a power cycle."

# Row 10
$ws.Cells.Item(10, 1).Value = 44331
$ws.Cells.Item(10, 2).Value = 44349
$ws.Cells.Item(10, 3).Value = 44349
$ws.Cells.Item(10, 4).Value = "PR-009"
$ws.Cells.Item(10, 5).Value = "Closed"
$ws.Cells.Item(10, 6).Value = "Software-app"
$ws.Cells.Item(10, 7).Value = "SS not aligned to SCK"
$ws.Cells.Item(10, 8).Value = "SS provides plenty of time ahead of SCK, but it doesn't deassert when SCK does."

# Row 11
$ws.Cells.Item(11, 1).Value = 44348
$ws.Cells.Item(11, 2).Value = 44349
$ws.Cells.Item(11, 3).Value = 44349
$ws.Cells.Item(11, 4).Value = "PR-010"
$ws.Cells.Item(11, 5).Value = "Closed"
$ws.Cells.Item(11, 6).Value = "FPGA app"
$ws.Cells.Item(11, 7).Value = "LD not functioning"
$ws.Cells.Item(11, 8).Value = "LD doesn't function, which prevents reads from occurring."

# Row 12
$ws.Cells.Item(12, 1).Value = 44379
$ws.Cells.Item(12, 4).Value = "PR-011"
$ws.Cells.Item(12, 5).Value = "Open"
$ws.Cells.Item(12, 6).Value = "FPGA app"
$ws.Cells.Item(12, 7).Value = "Error in fourth data exchange"
$ws.Cells.Item(12, 8).Value = "Incorrect value exchanged on SPI bus on fourth transaction"

# Row 13
$ws.Cells.Item(13, 1).Value = 44406
$ws.Cells.Item(13, 4).Value = "PR-012"
$ws.Cells.Item(13, 5).Value = "Open"
$ws.Cells.Item(13, 6).Value = "Software-app"
$ws.Cells.Item(13, 7).Value = "Synthesizer doesn't see data"
$ws.Cells.Item(13, 8).Value = "Data looks OK on the logic analyzer, but the synthesizer doesn't see the data."

# Row 14
$ws.Cells.Item(14, 1).Value = 44429
$ws.Cells.Item(14, 4).Value = "PR-013"
$ws.Cells.Item(14, 5).Value = "Open"
$ws.Cells.Item(14, 6).Value = "PCB"
$ws.Cells.Item(14, 7).Value = "FPGA doesn't drive MIDI out 0"
$ws.Cells.Item(14, 8).Value = "FPGA is able to drive output on pin 29, but not on pin 26"

# --- Apply date formats by copying existing styles (avoids duplicate style creation) ---

# Style s="1" (numFmtId 15, "d-mmm-yy") - reuse the pre-existing style from A2
$ws.Range("A2").Copy()
$ws.Range("A7:A14").PasteSpecial(-4122)
$ws.Range("B11:C11").PasteSpecial(-4122)

# Style s="2" (numFmtId 14, "mm-dd-yy") - new style, create once then copy to the rest
$ws.Range("B7").NumberFormat = "mm-dd-yy"
$ws.Range("B7").Copy()
$ws.Range("B7:C7").PasteSpecial(-4122)
$ws.Range("B8:C8").PasteSpecial(-4122)
$ws.Range("B10:C10").PasteSpecial(-4122)

# --- Column width adjustments (engine-quantized best effort to match target widths) ---
$ws.Columns.Item(6).ColumnWidth = 13.5
$ws.Columns.Item(7).ColumnWidth = 19.666666666666668

# --- Update selection to match final active cell ---
$ws.Range("G14").Select()
